$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data for columns F, G, H (same new value for every data row 2-7)
$F = 0.5
$G = 0.2098545
$H = 0.419709

# Row 2 (ECs)
$ws.Cells.Item(2,6).Value = $F
$ws.Cells.Item(2,7).Value = $G
$ws.Cells.Item(2,8).Value = $H
$ws.Cells.Item(2,13).Value = 30.7693535
$ws.Cells.Item(2,14).Value = 61.538707
$ws.Cells.Item(2,15).Value = 0.1179731387858698
$ws.Cells.Item(2,16).Value = 0.08351747770158975
$ws.Cells.Item(2,17).Value = 6.45708729406575
$ws.Cells.Item(2,18).Value = 25.828349176263
$ws.Cells.Item(2,19).Value = 0.1179731387858698
$ws.Cells.Item(2,20).Value = 0.08351747770158975

# Row 3 (FAPs)
$ws.Cells.Item(3,6).Value = $F
$ws.Cells.Item(3,7).Value = $G
$ws.Cells.Item(3,8).Value = $H
$ws.Cells.Item(3,15).Value = 0.09273042782012855
$ws.Cells.Item(3,16).Value = 0.09847086613229204
$ws.Cells.Item(3,17).Value = 5.075464410059
$ws.Cells.Item(3,18).Value = 30.452786460354
$ws.Cells.Item(3,19).Value = 0.09273042782012855
$ws.Cells.Item(3,20).Value = 0.09847086613229204

# Row 4 (Inflammatory-Mac)
$ws.Cells.Item(4,6).Value = $F
$ws.Cells.Item(4,7).Value = $G
$ws.Cells.Item(4,8).Value = $H
$ws.Cells.Item(4,13).Value = 81.77185533333333
$ws.Cells.Item(4,14).Value = 245.315566
$ws.Cells.Item(4,15).Value = 0.3135224286729781
$ws.Cells.Item(4,16).Value = 0.3329309033622996
$ws.Cells.Item(4,17).Value = 17.160191815049
$ws.Cells.Item(4,18).Value = 102.961150890294
$ws.Cells.Item(4,19).Value = 0.3135224286729781
$ws.Cells.Item(4,20).Value = 0.3329309033622996

# Row 5 (MuSCs)
$ws.Cells.Item(5,6).Value = $F
$ws.Cells.Item(5,7).Value = $G
$ws.Cells.Item(5,8).Value = $H
$ws.Cells.Item(5,13).Value = 14.8441875
$ws.Cells.Item(5,14).Value = 29.688375
$ws.Cells.Item(5,15).Value = 0.05691427322647431
$ws.Cells.Item(5,16).Value = 0.04029168498874919
$ws.Cells.Item(5,17).Value = 3.11511954571875
$ws.Cells.Item(5,18).Value = 12.460478182875
$ws.Cells.Item(5,19).Value = 0.05691427322647431
$ws.Cells.Item(5,20).Value = 0.04029168498874919

# Row 6 (Neutrophils)
$ws.Cells.Item(6,6).Value = $F
$ws.Cells.Item(6,7).Value = $G
$ws.Cells.Item(6,8).Value = $H
$ws.Cells.Item(6,13).Value = 53.27148833333334
$ws.Cells.Item(6,14).Value = 159.814465
$ws.Cells.Item(6,15).Value = 0.2042488376129897
$ws.Cells.Item(6,16).Value = 0.2168927763956593
$ws.Cells.Item(6,17).Value = 11.1792615484475
$ws.Cells.Item(6,18).Value = 67.07556929068501
$ws.Cells.Item(6,19).Value = 0.2042488376129897
$ws.Cells.Item(6,20).Value = 0.2168927763956593

# Row 7 (Resolving-Mac)
$ws.Cells.Item(7,6).Value = $F
$ws.Cells.Item(7,7).Value = $G
$ws.Cells.Item(7,8).Value = $H
$ws.Cells.Item(7,13).Value = 55.97408466666666
$ws.Cells.Item(7,14).Value = 167.922254
$ws.Cells.Item(7,15).Value = 0.2146108938815595
$ws.Cells.Item(7,16).Value = 0.22789629141941
$ws.Cells.Item(7,17).Value = 11.746413550681
$ws.Cells.Item(7,18).Value = 70.47848130408599
$ws.Cells.Item(7,19).Value = 0.2146108938815595
$ws.Cells.Item(7,20).Value = 0.22789629141941
